# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 10.611875054221343
$ws.Range("C2").Value = 20.493003485450572
$ws.Range("D2").Value = 25.139519920522048
$ws.Range("E2").Value = 25.85776053923928

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 13.102229730301303
$ws.Range("C3").Value = 16.817945928891334
$ws.Range("D3").Value = 13.129686290044475
$ws.Range("E3").Value = 29.383542863123807

# Update the sheet's active selection to match the narrower B1:E3 range
$ws.Range("B1:E3").Select()
